$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44547
$ws.Range("L2").Value = 'Especial'
$ws.Range("M2").Value = 350
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 20000
$ws.Range("R2").Value = 'Región Metropolitana'
$ws.Range("S2").Value = 1111

# Row 3
$ws.Range("D3").Value = 44547
$ws.Range("M3").Value = 350
$ws.Range("N3").Value = 18000
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 18000
$ws.Range("R3").Value = 'Región Metropolitana'
$ws.Range("S3").Value = 1000

# Row 4
$ws.Range("D4").Value = 44547
$ws.Range("M4").Value = 350
$ws.Range("N4").Value = 16000
$ws.Range("O4").Value = 16000
$ws.Range("P4").Value = 16000
$ws.Range("R4").Value = 'Región Metropolitana'
$ws.Range("S4").Value = 889

# Row 5
$ws.Range("D5").Value = 44187
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 350
$ws.Range("N5").Value = 16000
$ws.Range("P5").Value = 16000
$ws.Range("R5").Value = 'Región Metropolitana'
$ws.Range("S5").Value = 1067

# Row 6
$ws.Range("D6").Value = 44187
$ws.Range("K6").Value = 'Castle Brite'
$ws.Range("L6").Value = 'Segunda'
$ws.Range("M6").Value = 300
$ws.Range("N6").Value = 13000
$ws.Range("O6").Value = 13000
$ws.Range("P6").Value = 13000
$ws.Range("Q6").Value = '$/caja 15 kilos'
$ws.Range("S6").Value = 867
$ws.Range("T6").Value = 15

# Row 7
$ws.Range("L7").Value = 'Especial'
$ws.Range("N7").Value = 21000
$ws.Range("O7").Value = 21000
$ws.Range("P7").Value = 21000
$ws.Range("S7").Value = 1167

# Row 8
$ws.Range("L8").Value = 'Primera'
$ws.Range("N8").Value = 18000
$ws.Range("O8").Value = 18000
$ws.Range("P8").Value = 18000
$ws.Range("S8").Value = 1000

# Row 9
$ws.Range("D9").Value = 44575
$ws.Range("K9").Value = 'Modesto'
$ws.Range("M9").Value = 200
$ws.Range("O9").Value = 16000
$ws.Range("P9").Value = 16000
$ws.Range("Q9").Value = '$/caja 18 kilos'
$ws.Range("R9").Value = 'Región Metropolitana'
$ws.Range("S9").Value = 889
$ws.Range("T9").Value = 18

# Row 10
$ws.Range("D10").Value = 44162
$ws.Range("L10").Value = 'Tercera'
$ws.Range("N10").Value = 15000
$ws.Range("O10").Value = 16000
$ws.Range("P10").Value = 15500
$ws.Range("Q10").Value = '$/caja 15 kilos'
$ws.Range("R10").Value = 'Región de O''Higgins'
$ws.Range("S10").Value = 1033
$ws.Range("T10").Value = 15

# Row 11
$ws.Range("D11").Value = 44166
$ws.Range("L11").Value = 'Segunda'
$ws.Range("M11").Value = 600
$ws.Range("O11").Value = 17000
$ws.Range("P11").Value = 16500
$ws.Range("R11").Value = 'Región de O''Higgins'
$ws.Range("S11").Value = 1100

# Row 12
$ws.Range("D12").Value = 44530
$ws.Range("M12").Value = 500
$ws.Range("N12").Value = 20000
$ws.Range("O12").Value = 21000
$ws.Range("P12").Value = 20500
$ws.Range("Q12").Value = '$/caja 18 kilos'
$ws.Range("R12").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S12").Value = 1139
$ws.Range("T12").Value = 18

# Row 13
$ws.Range("D13").Value = 44553
$ws.Range("L13").Value = 'Especial'
$ws.Range("M13").Value = 250
$ws.Range("N13").Value = 20000
$ws.Range("O13").Value = 20000
$ws.Range("P13").Value = 20000
$ws.Range("Q13").Value = '$/caja 18 kilos'
$ws.Range("R13").Value = 'Región Metropolitana'
$ws.Range("S13").Value = 1111
$ws.Range("T13").Value = 18

# Row 14
$ws.Range("L14").Value = 'Primera'
$ws.Range("N14").Value = 18000
$ws.Range("O14").Value = 18000
$ws.Range("P14").Value = 18000
$ws.Range("S14").Value = 1000

# Row 15
$ws.Range("L15").Value = 'Segunda'
$ws.Range("N15").Value = 16000
$ws.Range("O15").Value = 16000
$ws.Range("P15").Value = 16000
$ws.Range("S15").Value = 889

# Row 16
$ws.Range("D16").Value = 44572
$ws.Range("K16").Value = 'Modesto'
$ws.Range("L16").Value = 'Especial'
$ws.Range("M16").Value = 150
$ws.Range("N16").Value = 21000
$ws.Range("O16").Value = 21000
$ws.Range("P16").Value = 21000
$ws.Range("S16").Value = 1167

# Row 17
$ws.Range("D17").Value = 44572
$ws.Range("K17").Value = 'Modesto'
$ws.Range("M17").Value = 150
$ws.Range("N17").Value = 18000
$ws.Range("O17").Value = 18000
$ws.Range("P17").Value = 18000
$ws.Range("R17").Value = 'Región Metropolitana'
$ws.Range("S17").Value = 1000

# Row 18
$ws.Range("D18").Value = 44572
$ws.Range("K18").Value = 'Modesto'
$ws.Range("M18").Value = 150
$ws.Range("N18").Value = 16000
$ws.Range("O18").Value = 16000
$ws.Range("P18").Value = 16000
$ws.Range("R18").Value = 'Región Metropolitana'
$ws.Range("S18").Value = 889

# Row 19
$ws.Range("D19").Value = 44176
$ws.Range("L19").Value = 'Segunda'
$ws.Range("M19").Value = 500
$ws.Range("N19").Value = 15000
$ws.Range("O19").Value = 16000
$ws.Range("P19").Value = 15500
$ws.Range("Q19").Value = '$/caja 15 kilos'
$ws.Range("R19").Value = 'Región Metropolitana'
$ws.Range("S19").Value = 1033
$ws.Range("T19").Value = 15

# Row 22
$ws.Range("D22").Value = 44537
$ws.Range("L22").Value = 'Primera'
$ws.Range("M22").Value = 500
$ws.Range("N22").Value = 20000
$ws.Range("O22").Value = 22000
$ws.Range("R22").Value = 'Provincia de San Felipe de Aconcagua'

# Row 23
$ws.Range("D23").Value = 44537
$ws.Range("L23").Value = 'Segunda'
$ws.Range("M23").Value = 250
$ws.Range("N23").Value = 17000
$ws.Range("O23").Value = 17000
$ws.Range("P23").Value = 17000
$ws.Range("R23").Value = 'Región del Maule'
$ws.Range("S23").Value = 944

# Row 24
$ws.Range("L24").Value = 'Especial'
$ws.Range("N24").Value = 21000
$ws.Range("O24").Value = 21000
$ws.Range("P24").Value = 21000
$ws.Range("S24").Value = 1167

# Row 25
$ws.Range("D25").Value = 44568
$ws.Range("L25").Value = 'Primera'
$ws.Range("M25").Value = 200
$ws.Range("N25").Value = 18000
$ws.Range("O25").Value = 18000
$ws.Range("P25").Value = 18000
$ws.Range("Q25").Value = '$/caja 18 kilos'
$ws.Range("S25").Value = 1000
$ws.Range("T25").Value = 18

# Row 26
$ws.Range("D26").Value = 44568
$ws.Range("K26").Value = 'Castle Brite'
$ws.Range("L26").Value = 'Segunda'
$ws.Range("N26").Value = 16000
$ws.Range("O26").Value = 16000
$ws.Range("P26").Value = 16000
$ws.Range("S26").Value = 889

# Row 27
$ws.Range("D27").Value = 44169
$ws.Range("K27").Value = 'Castle Brite'
$ws.Range("L27").Value = 'Segunda'
$ws.Range("M27").Value = 500
$ws.Range("N27").Value = 15000
$ws.Range("O27").Value = 16000
$ws.Range("P27").Value = 15500
$ws.Range("Q27").Value = '$/caja 15 kilos'
$ws.Range("R27").Value = 'Región de O''Higgins'
$ws.Range("S27").Value = 1033
$ws.Range("T27").Value = 15

# Row 28
$ws.Range("D28").Value = 44540
$ws.Range("K28").Value = 'Castle Brite'
$ws.Range("M28").Value = 600
$ws.Range("R28").Value = 'Región del Maule'

# Row 29
$ws.Range("D29").Value = 44194
$ws.Range("M29").Value = 300

# Row 33
$ws.Range("D33").Value = 44551
$ws.Range("L33").Value = 'Especial'
$ws.Range("M33").Value = 200
$ws.Range("N33").Value = 20000
$ws.Range("O33").Value = 20000
$ws.Range("P33").Value = 20000
$ws.Range("Q33").Value = '$/caja 18 kilos'
$ws.Range("R33").Value = 'Región Metropolitana'
$ws.Range("S33").Value = 1111
$ws.Range("T33").Value = 18

# Row 34
$ws.Range("L34").Value = 'Primera'
$ws.Range("N34").Value = 18000
$ws.Range("O34").Value = 18000
$ws.Range("P34").Value = 18000
$ws.Range("S34").Value = 1000

# Row 35
$ws.Range("L35").Value = 'Segunda'
$ws.Range("N35").Value = 16000
$ws.Range("O35").Value = 16000
$ws.Range("P35").Value = 16000
$ws.Range("S35").Value = 889

# Row 36
$ws.Range("D36").Value = 44579
$ws.Range("K36").Value = 'Modesto'
$ws.Range("L36").Value = 'Especial'
$ws.Range("N36").Value = 21000
$ws.Range("O36").Value = 21000
$ws.Range("P36").Value = 21000
$ws.Range("S36").Value = 1167

# Row 37
$ws.Range("D37").Value = 44579
$ws.Range("K37").Value = 'Modesto'
$ws.Range("L37").Value = 'Primera'
$ws.Range("M37").Value = 200
$ws.Range("N37").Value = 18000
$ws.Range("O37").Value = 18000
$ws.Range("P37").Value = 18000
$ws.Range("S37").Value = 1000

# Row 38
$ws.Range("D38").Value = 44579
$ws.Range("K38").Value = 'Modesto'
$ws.Range("L38").Value = 'Segunda'
$ws.Range("M38").Value = 200
$ws.Range("N38").Value = 16000
$ws.Range("O38").Value = 16000
$ws.Range("P38").Value = 16000
$ws.Range("S38").Value = 889

# Row 39
$ws.Range("D39").Value = 44533
$ws.Range("L39").Value = 'Primera'
$ws.Range("N39").Value = 24000
$ws.Range("O39").Value = 24000
$ws.Range("P39").Value = 24000
$ws.Range("R39").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S39").Value = 1333

# Row 40
$ws.Range("D40").Value = 44533
$ws.Range("K40").Value = 'Castle Brite'
$ws.Range("L40").Value = 'Segunda'
$ws.Range("M40").Value = 350
$ws.Range("N40").Value = 20000
$ws.Range("O40").Value = 20000
$ws.Range("P40").Value = 20000
$ws.Range("R40").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S40").Value = 1111

# Row 41
$ws.Range("D41").Value = 44533
$ws.Range("K41").Value = 'Castle Brite'
$ws.Range("L41").Value = 'Tercera'
$ws.Range("M41").Value = 350
$ws.Range("N41").Value = 17000
$ws.Range("O41").Value = 17000
$ws.Range("P41").Value = 17000
$ws.Range("R41").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S41").Value = 944

# Row 42
$ws.Range("D42").Value = 44159
$ws.Range("K42").Value = 'Castle Brite'
$ws.Range("L42").Value = 'Tercera'
$ws.Range("M42").Value = 400
$ws.Range("N42").Value = 15500
$ws.Range("P42").Value = 15750
$ws.Range("Q42").Value = '$/caja 15 kilos'
$ws.Range("R42").Value = 'Región de O''Higgins'
$ws.Range("S42").Value = 1050
$ws.Range("T42").Value = 15
